$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (Coin names, Links, double-dot prices, and Volume % strings)
# These are safe to assign directly without numeric auto-conversion.
$ws.Range('D2').Value = '62.520.39'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').Value = '3.436.18'
$ws.Range('E3').Value = '  +2.46%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('E6').Value = '  +1.43%  '
$ws.Range('E7').Value = '  -0.89%  '
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('E9').Value = '  +3.46%  '
$ws.Range('E10').Value = '  +9.10%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E14').Value = '  +1.44%  '
$ws.Range('D15').Value = '3.484.98'
$ws.Range('E15').Value = '  +4.19%  '
$ws.Range('D16').Value = '62.470.93'
$ws.Range('E16').Value = '  +1.69%  '
$ws.Range('E17').Value = '  +1.71%  '
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('E19').Value = '  +16.23%  '
$ws.Range('E20').Value = '  -1.84%  '
$ws.Range('E21').Value = '  +3.18%  '
$ws.Range('E22').Value = '  +2.74%  '
$ws.Range('E23').Value = '  -1.55%  '
$ws.Range('E24').Value = '  +1.06%  '
$ws.Range('E25').Value = '  +1.68%  '
$ws.Range('E26').Value = '  +1.08%  '
$ws.Range('E27').Value = '  -4.06%  '
$ws.Range('E28').Value = '  +5.27%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E29').Value = '  +7.56%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E30').Value = '  +6.40%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('E33').Value = '  -2.24%  '
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  -1.00%  '
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('E39').Value = '  +13.76%  '
$ws.Range('E40').Value = '  -1.75%  '
$ws.Range('E41').Value = '  +4.61%  '
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('E46').Value = '  -0.50%  '
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').Value = '2.112.76'
$ws.Range('E48').Value = '  -1.11%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('E49').Value = '  -1.69%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('E50').Value = '  +3.99%  '
$ws.Range('E51').Value = '  +28.36%  '

# Price values that look like pure decimal numbers must be forced to Text
# so Excel keeps exact digits (e.g. trailing zeros, no float rounding),
# matching the original inline-string cell formatting.
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '407.45'
$cell.Style = 'Normal'
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '130.24'
$cell.Style = 'Normal'
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.595'
$cell.Style = 'Normal'
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.691'
$cell.Style = 'Normal'
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.136'
$cell.Style = 'Normal'
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '41.95'
$cell.Style = 'Normal'
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '8.43'
$cell.Style = 'Normal'
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '19.81'
$cell.Style = 'Normal'
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '11.59'
$cell.Style = 'Normal'
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.0000153'
$cell.Style = 'Normal'
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '84.52'
$cell.Style = 'Normal'
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '313.03'
$cell.Style = 'Normal'
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '3.17'
$cell.Style = 'Normal'
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '29.73'
$cell.Style = 'Normal'
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '8.13'
$cell.Style = 'Normal'
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '7.80'
$cell.Style = 'Normal'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '2.79'
$cell.Style = 'Normal'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '44.49'
$cell.Style = 'Normal'
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '51.78'
$cell.Style = 'Normal'
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '2.97'
$cell.Style = 'Normal'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.321'
$cell.Style = 'Normal'
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '3.33'
$cell.Style = 'Normal'
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '141.81'
$cell.Style = 'Normal'
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.125'
$cell.Style = 'Normal'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '3.92'
$cell.Style = 'Normal'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '16.83'
$cell.Style = 'Normal'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '21.45'
$cell.Style = 'Normal'
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '2.32'
$cell.Style = 'Normal'
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '1.96'
$cell.Style = 'Normal'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '1.09'
$cell.Style = 'Normal'
